$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fix row 10: I10/J10 reset to 0 (bug fix referenced in commit message) ---
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0

# --- New row 42: 2NH01 combo ---
$ws.Range("A42").Value = "2NH01"
$ws.Range("B42").Value = "COMBO 2 NƯỚC HOA NHÀI"
$ws.Range("C42").Value = 1
$ws.Range("D42").Value = "NH01"
$ws.Range("E42").Value = "NH01"
$ws.Range("F42").Value = "Nước hoa nhài"
$ws.Range("G42").Value = 2
$ws.Range("H42").Value = "Lọ"
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0

# --- New row 43: 3NH01 combo ---
$ws.Range("A43").Value = "3NH01"
$ws.Range("B43").Value = "COMBO 3 NƯỚC HOA NHÀI"
$ws.Range("C43").Value = 1
$ws.Range("D43").Value = "NH01"
$ws.Range("E43").Value = "NH01"
$ws.Range("F43").Value = "Nước hoa nhài"
$ws.Range("G43").Value = 3
$ws.Range("H43").Value = "Lọ"
$ws.Range("I43").Value = 140909
$ws.Range("J43").Formula = "=I43*3"

# --- New row 44: MAS1 Massage machine ---
$ws.Range("A44").Value = "MAS1"
$ws.Range("B44").Value = "Máy Massage Cổ vai gáy"
$ws.Range("C44").Value = 1
$ws.Range("D44").Value = "MAS1"
$ws.Range("E44").Value = "MAS1"
$ws.Range("F44").Value = "Máy massage cổ dùng pin sạc 6D"
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = "Cái"
$ws.Range("I44").Value = 1150000
$ws.Range("J44").Value = 1150000

# --- Font styling: product-name cells pasted with Segoe UI 7pt formatting ---
$ws.Range("B42").Font.Name = "Segoe UI"
$ws.Range("B42").Font.Size = 7
$ws.Range("B43").Font.Name = "Segoe UI"
$ws.Range("B43").Font.Size = 7
$ws.Range("B44").Font.Name = "Segoe UI"
$ws.Range("B44").Font.Size = 7
$ws.Range("F44").Font.Name = "Segoe UI"
$ws.Range("F44").Font.Size = 7

# --- Hyperlinks for the SKU codes linking back to the storefront admin ---
$url1 = "https://caoxoaquoccoquocnghiep.mysapogo.com/admin/products/193913949/variants/315993725"
$url2 = "https://caoxoaquoccoquocnghiep.mysapogo.com/admin/products/193913988/variants/315993781"
$ws.Hyperlinks.Add($ws.Range("A42"), $url1, "", "", $url1)
$ws.Range("A42").Value = "2NH01"
$ws.Hyperlinks.Add($ws.Range("A43"), $url2, "", "", $url2)
$ws.Range("A43").Value = "3NH01"

# --- Restore the view to where the user left off editing ---
[void]$ws.Range("F41").Select()
